$d = $word.ActiveDocument

# 1. Fix the title wording: "DESCARGUIO" -> "DESCARGO"
#    (Match case so we don't touch the unrelated lowercase
#    "...orden de descarguio..." sentence elsewhere in the document.)
$d.Content.Find.Execute("DESCARGUIO", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "DESCARGO", 2) | Out-Null

# 2. Clear the leftover "no fill" shading on the template data row of the
#    second table (the header row keeps its teal ADFFFF shading; only the
#    data row underneath, whose shading was just the default/auto, loses
#    its explicit shading formatting).
$t = $d.Tables.Item(2)
$row = $t.Rows.Item(2)
for ($c = 1; $c -le $row.Cells.Count; $c++) {
    $cell = $row.Cells.Item($c)
    $shd = $cell.Shading
    $shd.Texture = 0
    $shd.ForegroundPatternColor = -16777216
    $shd.BackgroundPatternColor = -16777216
}
